$wb = $excel.ActiveWorkbook

# ---------- Sheet ALC ----------
$ws = $wb.Worksheets.Item("ALC")

# Row 125
$ws.Range("H125").Value = 999
$ws.Range("I125").Value = 998
$ws.Range("J125").Value = 999.5
$ws.Range("K125").Value = 8982
$ws.Range("L125").Value = 8995.5
$ws.Range("M125").Value = -6522
$ws.Range("N125").Value = -13915.5

# Row 137
$ws.Range("H137").Value = 2354.111
$ws.Range("I137").Value = 1915.5385
$ws.Range("J137").Value = 2761.3572
$ws.Range("K137").Value = 5746.6155
$ws.Range("L137").Value = 8284.071599999999
$ws.Range("M137").Value = -3196.6155
$ws.Range("N137").Value = -13384.0716

# ---------- Sheet ARM ----------
$ws = $wb.Worksheets.Item("ARM")

# Row 63
$ws.Range("H63").Value = 13776.667
$ws.Range("J63").Value = 13776.667
$ws.Range("L63").Value = 13776.667
$ws.Range("N63").Value = -15148.667

# Row 66
$ws.Range("H66").Value = 13776.667
$ws.Range("J66").Value = 13776.667
$ws.Range("L66").Value = 68883.33499999999
$ws.Range("N66").Value = -75747.33499999999

# Row 88
$ws.Range("H88").Value = 3701.2
$ws.Range("J88").Value = 5166.6665
$ws.Range("L88").Value = 5166.6665
$ws.Range("N88").Value = -5978.6665

# Row 91
$ws.Range("H91").Value = 3701.2
$ws.Range("J91").Value = 5166.6665
$ws.Range("L91").Value = 5166.6665
$ws.Range("N91").Value = -7974.6665

# Row 109
$ws.Range("H109").Value = 87500
$ws.Range("J109").Value = 87500
$ws.Range("L109").Value = 87500
$ws.Range("N109").Value = -90274

# ---------- Sheet BSM ----------
$ws = $wb.Worksheets.Item("BSM")

# Row 105
$ws.Range("H105").Value = 10139.2

# ---------- Sheet CRP ----------
$ws = $wb.Worksheets.Item("CRP")

# Row 74
$ws.Range("H74").Value = 59984.668
$ws.Range("J74").Value = 59984.668
$ws.Range("L74").Value = 59984.668
$ws.Range("N74").Value = -61732.668

# Row 77
$ws.Range("H77").Value = 59984.668
$ws.Range("J77").Value = 59984.668
$ws.Range("L77").Value = 179954.004
$ws.Range("N77").Value = -188690.004

# Row 122
$ws.Range("H122").Value = 1795
$ws.Range("I122").Value = 1795
$ws.Range("K122").Value = 5385
$ws.Range("M122").Value = -2935

# Row 132
$ws.Range("H132").Value = 2474.75
$ws.Range("I132").Value = 1260
$ws.Range("K132").Value = 3780
$ws.Range("M132").Value = -1250

# ---------- Sheet CUL ----------
$ws = $wb.Worksheets.Item("CUL")

# Row 107
$ws.Range("H107").Value = 545.2143
$ws.Range("I107").Value = 357.14285
$ws.Range("K107").Value = 1071.42855
$ws.Range("M107").Value = 848.5714499999999

# Row 137
$ws.Range("H137").Value = 5199.75
$ws.Range("I137").Value = 5199.75
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 15599.25
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -10499.25
$ws.Range("N137").ClearContents()

# Row 139
$ws.Range("H139").Value = 128037.5
$ws.Range("I139").Value = 334766.66
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 1004299.98
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -999159.98
$ws.Range("N139").Value = -22280

# Row 141
$ws.Range("H141").Value = 27600
$ws.Range("I141").Value = 27600
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 82800
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -77620
$ws.Range("N141").ClearContents()

# ---------- Sheet GSM ----------
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 9110.1
$ws.Range("J80").Value = 20666
$ws.Range("L80").Value = 20666
$ws.Range("N80").Value = -22662

# Row 83
$ws.Range("H83").Value = 9110.1
$ws.Range("J83").Value = 20666
$ws.Range("L83").Value = 103330
$ws.Range("N83").Value = -113314

# Rows 125-141: clear columns H through N entirely
for ($r = 125; $r -le 141; $r++) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---------- Sheet LTW ----------
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 2285
$ws.Range("I61").Value = 2299
$ws.Range("K61").Value = 2299
$ws.Range("M61").Value = -2097

# Row 113
$ws.Range("H113").Value = 2285
$ws.Range("I113").Value = 2299
$ws.Range("K113").Value = 2299
$ws.Range("M113").Value = -129
